$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2022.5
$ws.Range("J17").Value = 2022.5
$ws.Range("L17").Value = 6067.5
$ws.Range("N17").Value = -6403.5
$ws.Range("H28").Value = 2225.3
$ws.Range("I28").Value = 531.6875
$ws.Range("J28").Value = 8999.75
$ws.Range("K28").Value = 531.6875
$ws.Range("L28").Value = 8999.75
$ws.Range("M28").Value = -46.6875
$ws.Range("N28").Value = -9969.75
$ws.Range("H64").Value = 83340490
$ws.Range("I64").Value = 8109
$ws.Range("J64").Value = 125006680
$ws.Range("K64").Value = 8109
$ws.Range("L64").Value = 125006680
$ws.Range("M64").Value = -7861
$ws.Range("N64").Value = -125007176
$ws.Range("H67").Value = 83340490
$ws.Range("I67").Value = 8109
$ws.Range("J67").Value = 125006680
$ws.Range("K67").Value = 8109
$ws.Range("L67").Value = 125006680
$ws.Range("M67").Value = -7251
$ws.Range("N67").Value = -125008396
$ws.Range("H74").Value = 15128.84
$ws.Range("I74").Value = 15764.263
$ws.Range("K74").Value = 15764.263
$ws.Range("M74").Value = -14828.263
$ws.Range("H77").Value = 15128.84
$ws.Range("I77").Value = 15764.263
$ws.Range("K77").Value = 78821.315
$ws.Range("M77").Value = -74141.315
$ws.Range("H116").Value = 5083
$ws.Range("I116").Value = 3124.5
$ws.Range("K116").Value = 3124.5
$ws.Range("M116").Value = 317.5
$ws.Range("H137").Value = 1435456.2
$ws.Range("J137").Value = 7868.9
$ws.Range("L137").Value = 23606.7
$ws.Range("N137").Value = -28706.7
$ws.Range("H141").Value = 3811.9092
$ws.Range("I141").Value = 3669.5557
$ws.Range("K141").Value = 11008.6671
$ws.Range("M141").Value = -5828.667099999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("H32").Value = 3519.3704
$ws.Range("I32").Value = 3022.8696
$ws.Range("J32").Value = 6374.25
$ws.Range("K32").Value = 3022.8696
$ws.Range("L32").Value = 6374.25
$ws.Range("M32").Value = -2735.8696
$ws.Range("N32").Value = -6948.25
$ws.Range("H45").Value = 34076.848
$ws.Range("I45").Value = 39636.09
$ws.Range("K45").Value = 39636.09
$ws.Range("M45").Value = -39259.09
$ws.Range("H61").Value = 5642
$ws.Range("I61").Value = 2989.5
$ws.Range("K61").Value = 2989.5
$ws.Range("M61").Value = -2777.5
$ws.Range("H97").Value = 1280.4073
$ws.Range("I97").Value = 1152.8
$ws.Range("J97").Value = 1645
$ws.Range("K97").Value = 1152.8
$ws.Range("L97").Value = 1645
$ws.Range("M97").Value = -656.8
$ws.Range("N97").Value = -2637
$ws.Range("H102").Value = 3787.818
$ws.Range("I102").Value = 3111
$ws.Range("K102").Value = 3111
$ws.Range("M102").Value = -1489
$ws.Range("H132").Value = 3064.6667
$ws.Range("I132").Value = 1972.125
$ws.Range("K132").Value = 5916.375
$ws.Range("M132").Value = -3386.375
$ws.Range("H136").Value = 5642
$ws.Range("I136").Value = 2989.5
$ws.Range("K136").Value = 8968.5
$ws.Range("M136").Value = -6418.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1448.92
$ws.Range("I107").Value = 1287.7646
$ws.Range("J107").Value = 1791.375
$ws.Range("K107").Value = 1287.7646
$ws.Range("L107").Value = 1791.375
$ws.Range("M107").Value = 632.2354
$ws.Range("N107").Value = -5631.375
$ws.Range("H134").Value = 4566.067
$ws.Range("I134").Value = 4499.25
$ws.Range("J134").Value = 4833.3335
$ws.Range("K134").Value = 13497.75
$ws.Range("L134").Value = 14500.0005
$ws.Range("M134").Value = -10962.75
$ws.Range("N134").Value = -19570.0005
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1379.6471
$ws.Range("I105").Value = 832.75
$ws.Range("J105").Value = 1865.7778
$ws.Range("K105").Value = 832.75
$ws.Range("L105").Value = 1865.7778
$ws.Range("M105").Value = 914.25
$ws.Range("N105").Value = -5359.7778
$ws.Range("H132").Value = 3377.5
$ws.Range("I132").Value = 3377.5
$ws.Range("K132").Value = 10132.5
$ws.Range("M132").Value = -7602.5
$ws.Range("H134").Value = 3049
$ws.Range("I134").Value = 2943.3333
$ws.Range("K134").Value = 8829.999899999999
$ws.Range("M134").Value = -6294.999899999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 30567750
$ws.Range("J4").Value = 6243193.5
$ws.Range("L4").Value = 18729580.5
$ws.Range("N4").Value = -18729804.5
$ws.Range("H6").Value = 150
$ws.Range("I6").Value = 150
$ws.Range("K6").Value = 450
$ws.Range("M6").Value = -337
$ws.Range("H9").Value = 276033.34
$ws.Range("I9").Value = 276033.34
$ws.Range("K9").Value = 828100.02
$ws.Range("M9").Value = -827876.02
$ws.Range("H11").Value = 933.2778
$ws.Range("I11").Value = 964.64703
$ws.Range("K11").Value = 2893.94109
$ws.Range("M11").Value = -2753.94109
$ws.Range("H13").Value = 500
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 500
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 1500
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -1836
$ws.Range("H15").Value = 936.375
$ws.Range("I15").Value = 81.833336
$ws.Range("J15").Value = 3500
$ws.Range("K15").Value = 245.500008
$ws.Range("L15").Value = 10500
$ws.Range("M15").Value = -105.500008
$ws.Range("N15").Value = -10780
$ws.Range("H16").Value = 1166.3334
$ws.Range("I16").Value = 374.5
$ws.Range("J16").Value = 2750
$ws.Range("K16").Value = 1123.5
$ws.Range("L16").Value = 8250
$ws.Range("M16").Value = -950.5
$ws.Range("N16").Value = -8596
$ws.Range("H17").Value = 153.33333
$ws.Range("I17").Value = 153.33333
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 459.99999
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -290.99999
$ws.Range("N17").ClearContents()
$ws.Range("H60").Value = 1056032.8
$ws.Range("I60").Value = 3333723.8
$ws.Range("K60").Value = 10001171.4
$ws.Range("M60").Value = -10000920.4
$ws.Range("H121").Value = 333766.66
$ws.Range("I121").Value = 300
$ws.Range("J121").Value = 500500
$ws.Range("K121").Value = 900
$ws.Range("L121").Value = 1501500
$ws.Range("M121").Value = 410
$ws.Range("N121").Value = -1504120
$ws.Range("H125").Value = 7000
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H131").Value = 1932.1897
$ws.Range("I131").Value = 2432
$ws.Range("J131").Value = 1852.22
$ws.Range("K131").Value = 7296
$ws.Range("L131").Value = 5556.66
$ws.Range("M131").Value = -2256
$ws.Range("N131").Value = -15636.66
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 14998.75
$ws.Range("I12").Value = 25000
$ws.Range("J12").Value = 4997.5
$ws.Range("K12").Value = 25000
$ws.Range("L12").Value = 4997.5
$ws.Range("M12").Value = -24830
$ws.Range("N12").Value = -5337.5
$ws.Range("H61").Value = 9562.786
$ws.Range("I61").Value = 1943.091
$ws.Range("J61").Value = 37501.668
$ws.Range("K61").Value = 1943.091
$ws.Range("L61").Value = 37501.668
$ws.Range("M61").Value = -1741.091
$ws.Range("N61").Value = -37905.668
$ws.Range("H113").Value = 9562.786
$ws.Range("I113").Value = 1943.091
$ws.Range("J113").Value = 37501.668
$ws.Range("K113").Value = 1943.091
$ws.Range("L113").Value = 37501.668
$ws.Range("M113").Value = 226.9090000000001
$ws.Range("N113").Value = -41841.668
$ws.Range("H135").Value = 32998
$ws.Range("J135").Value = 32998
$ws.Range("L135").Value = 32998
$ws.Range("N135").Value = -43138
$ws.Range("H136").Value = 4811.1333
$ws.Range("I136").Value = 5541.1113
$ws.Range("J136").Value = 3716.1667
$ws.Range("K136").Value = 16623.3339
$ws.Range("L136").Value = 11148.5001
$ws.Range("M136").Value = -14073.3339
$ws.Range("N136").Value = -16248.5001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 40002.5
$ws.Range("I13").Value = 40005
$ws.Range("J13").Value = 40000
$ws.Range("K13").Value = 40005
$ws.Range("L13").Value = 40000
$ws.Range("M13").Value = -39865
$ws.Range("N13").Value = -40280
